$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1232, 1).Value = "U42_01"
$ws.Cells.Item(1232, 2).Value = 42
$ws.Cells.Item(1232, 3).Value = "Tai nghe"
$ws.Cells.Item(1232, 4).Value = "Headphones"
$ws.Cells.Item(1232, 5).Value = "Use a set of headphones to listen clearly"
$ws.Cells.Item(1232, 6).Value = "set of headphones / bộ tai nghe"
$ws.Cells.Item(1232, 7).Value = "N"

$ws.Cells.Item(1233, 1).Value = "U42_02"
$ws.Cells.Item(1233, 2).Value = 42
$ws.Cells.Item(1233, 3).Value = "Cổ điển"
$ws.Cells.Item(1233, 4).Value = "Classical"
$ws.Cells.Item(1233, 5).Value = "The classical concert is at the HaNoi Opera House"
$ws.Cells.Item(1233, 6).Value = "A classical concert / buổi hòa nhạc cổ điển"
$ws.Cells.Item(1233, 7).Value = "Adj"

$ws.Cells.Item(1234, 1).Value = "U42_03"
$ws.Cells.Item(1234, 2).Value = 42
$ws.Cells.Item(1234, 3).Value = "Nhạc cụ"
$ws.Cells.Item(1234, 4).Value = "Instrument"
$ws.Cells.Item(1234, 5).Value = "Can you play any instrument?"
$ws.Cells.Item(1234, 6).Value = "play an instrument"
$ws.Cells.Item(1234, 7).Value = "N"

$ws.Cells.Item(1235, 1).Value = "U42_04"
$ws.Cells.Item(1235, 2).Value = 42
$ws.Cells.Item(1235, 3).Value = "Cuộc thi"
$ws.Cells.Item(1235, 4).Value = "Contest"
$ws.Cells.Item(1235, 5).Value = "The contest between Jack and Jill is exciting."
$ws.Cells.Item(1235, 6).Value = "contest between A and B"
$ws.Cells.Item(1235, 7).Value = "N"

$ws.Cells.Item(1236, 1).Value = "U42_05"
$ws.Cells.Item(1236, 2).Value = 42
$ws.Cells.Item(1236, 3).Value = "Sáo"
$ws.Cells.Item(1236, 4).Value = "Flute"
$ws.Cells.Item(1236, 5).Value = "I am learning to play the flute"
$ws.Cells.Item(1236, 6).Value = "play the flute"
$ws.Cells.Item(1236, 7).Value = "N"

$ws.Cells.Item(1237, 1).Value = "U42_06"
$ws.Cells.Item(1237, 2).Value = 42
$ws.Cells.Item(1237, 3).Value = "Nhạc sĩ, nhạc công"
$ws.Cells.Item(1237, 4).Value = "Musicican"
$ws.Cells.Item(1237, 5).Value = "Professional musicians have to train very hard."
$ws.Cells.Item(1237, 6).Value = "a professional musician / nhạc công chuyên nghiệp"
$ws.Cells.Item(1237, 7).Value = "N"

$ws.Cells.Item(1238, 1).Value = "U42_07"
$ws.Cells.Item(1238, 2).Value = 42
$ws.Cells.Item(1238, 3).Value = "Thu âm"
$ws.Cells.Item(1238, 4).Value = "Record"
$ws.Cells.Item(1238, 5).Value = "come to the studio to record music"
$ws.Cells.Item(1238, 6).Value = "record something"
$ws.Cells.Item(1238, 7).Value = "V"

$ws.Cells.Item(1239, 1).Value = "U42_08"
$ws.Cells.Item(1239, 2).Value = 42
$ws.Cells.Item(1239, 3).Value = "hát theo"
$ws.Cells.Item(1239, 4).Value = "Sing along"
$ws.Cells.Item(1239, 5).Value = "if you know this song, please sing along with us"
$ws.Cells.Item(1239, 6).Value = "sing along with somebody"
$ws.Cells.Item(1239, 7).Value = "V"

$ws.Cells.Item(1240, 1).Value = "U42_09"
$ws.Cells.Item(1240, 2).Value = 42
$ws.Cells.Item(1240, 3).Value = "Bảng sếp hạng, bảng biểu"
$ws.Cells.Item(1240, 4).Value = "Chart"
$ws.Cells.Item(1240, 5).Value = "This song is number one in music charts these days"
$ws.Cells.Item(1240, 6).Value = "The music charts / bảng sếp hạng âm nhạc"
$ws.Cells.Item(1240, 7).Value = "N"

$ws.Cells.Item(1241, 1).Value = "U42_10"
$ws.Cells.Item(1241, 2).Value = 42
$ws.Cells.Item(1241, 3).Value = "Tông giọng"
$ws.Cells.Item(1241, 4).Value = "Tone"
$ws.Cells.Item(1241, 5).Value = "You can guess her feelings through the tone of her voice"
$ws.Cells.Item(1241, 6).Value = "tone of one's voice / tông giọng của một người"
$ws.Cells.Item(1241, 7).Value = "N"

$ws.Cells.Item(1242, 1).Value = "U42_11"
$ws.Cells.Item(1242, 2).Value = 42
$ws.Cells.Item(1242, 3).Value = "Giai điệu"
$ws.Cells.Item(1242, 4).Value = "Tune"
$ws.Cells.Item(1242, 5).Value = "I can dance to the tune of any song."
$ws.Cells.Item(1242, 6).Value = "to the tune of something / theo giai điệu của một cái gì đó"
$ws.Cells.Item(1242, 7).Value = "N"

$ws.Cells.Item(1243, 1).Value = "U42_12"
$ws.Cells.Item(1243, 2).Value = 42
$ws.Cells.Item(1243, 3).Value = "Buổi triển lãm"
$ws.Cells.Item(1243, 4).Value = "Exhibition"
$ws.Cells.Item(1243, 5).Value = "This is an exhibition of children's paintings."
$ws.Cells.Item(1243, 6).Value = "an exhibition of something"
$ws.Cells.Item(1243, 7).Value = "N"

$ws.Cells.Item(1244, 1).Value = "U42_13"
$ws.Cells.Item(1244, 2).Value = 42
$ws.Cells.Item(1244, 3).Value = "Nghệ sĩ"
$ws.Cells.Item(1244, 4).Value = "Artist"
$ws.Cells.Item(1244, 5).Value = "She works as a solo artist"
$ws.Cells.Item(1244, 6).Value = "a solo artist / nghệ sĩ độc tấu"
$ws.Cells.Item(1244, 7).Value = "N"

$ws.Cells.Item(1245, 1).Value = "U42_14"
$ws.Cells.Item(1245, 2).Value = 42
$ws.Cells.Item(1245, 3).Value = "Họa sĩ"
$ws.Cells.Item(1245, 4).Value = "Painter"
$ws.Cells.Item(1245, 5).Value = "Landscape painters paint natural scenery (vẽ cảnh quan thiên nhiên)"
$ws.Cells.Item(1245, 6).Value = "a landscape painter / họa sĩ tranh phong cảnh"
$ws.Cells.Item(1245, 7).Value = "N"

$ws.Cells.Item(1246, 1).Value = "U42_15"
$ws.Cells.Item(1246, 2).Value = 42
$ws.Cells.Item(1246, 3).Value = "Tác phẩm điêu khắc"
$ws.Cells.Item(1246, 4).Value = "Sculpture"
$ws.Cells.Item(1246, 5).Value = "The museum is full of marble sculpture"
$ws.Cells.Item(1246, 6).Value = "a marble sculpture / một tác phẩm điêu khắc bằng đá cẩm thạch"
$ws.Cells.Item(1246, 7).Value = "N"

$ws.Cells.Item(1247, 1).Value = "U42_16"
$ws.Cells.Item(1247, 2).Value = 42
$ws.Cells.Item(1247, 3).Value = "Tượng"
$ws.Cells.Item(1247, 4).Value = "Statue"
$ws.Cells.Item(1247, 5).Value = "The bronze statues need polishing ( lau chùi)"
$ws.Cells.Item(1247, 6).Value = "a bronze statue / bức tượng đồng"
$ws.Cells.Item(1247, 7).Value = "N"

$ws.Cells.Item(1248, 1).Value = "U42_17"
$ws.Cells.Item(1248, 2).Value = 42
$ws.Cells.Item(1248, 3).Value = "Trống"
$ws.Cells.Item(1248, 4).Value = "Drum"
$ws.Cells.Item(1248, 5).Value = "Can the musician play the drum? "
$ws.Cells.Item(1248, 6).Value = "Play the drum"
$ws.Cells.Item(1248, 7).Value = "N"

$ws.Cells.Item(1249, 1).Value = "U42_18"
$ws.Cells.Item(1249, 2).Value = 42
$ws.Cells.Item(1249, 3).Value = "Hòa trộn"
$ws.Cells.Item(1249, 4).Value = "Blend"
$ws.Cells.Item(1249, 5).Value = "This song is a great blend of modern and traditional music."
$ws.Cells.Item(1249, 6).Value = "a blend of / một sự pha trộn của"
$ws.Cells.Item(1249, 7).Value = "N"

$ws.Cells.Item(1250, 1).Value = "U42_19"
$ws.Cells.Item(1250, 2).Value = 42
$ws.Cells.Item(1250, 3).Value = "Thuộc dân gian"
$ws.Cells.Item(1250, 4).Value = "Folk"
$ws.Cells.Item(1250, 5).Value = "We learned many folk songs in kindergarten (trường mẫu giáo)"
$ws.Cells.Item(1250, 6).Value = "folk song / bài dân ca"
$ws.Cells.Item(1250, 7).Value = "N"

$ws.Cells.Item(1251, 1).Value = "U42_20"
$ws.Cells.Item(1251, 2).Value = 42
$ws.Cells.Item(1251, 3).Value = "Nhịp điệu"
$ws.Cells.Item(1251, 4).Value = "Rhythm"
$ws.Cells.Item(1251, 5).Value = "Please dance in rhythm with the music"
$ws.Cells.Item(1251, 6).Value = "in rhythm with something / đúng nhịp với cái gì"
$ws.Cells.Item(1251, 7).Value = "N"

$ws.Cells.Item(1252, 1).Value = "U42_21"
$ws.Cells.Item(1252, 2).Value = 42
$ws.Cells.Item(1252, 3).Value = "Hoàn toàn"
$ws.Cells.Item(1252, 4).Value = "absolutely"
$ws.Cells.Item(1252, 5).Value = "Absolutely nothing can go wrong now (có thể sai sót được nữa)"
$ws.Cells.Item(1252, 6).Value = "Absolutely nothing / hoàn toàn không có gì"
$ws.Cells.Item(1252, 7).Value = "Adj"

$ws.Cells.Item(1253, 1).Value = "U42_22"
$ws.Cells.Item(1253, 2).Value = 42
$ws.Cells.Item(1253, 3).Value = "Thơ ca"
$ws.Cells.Item(1253, 4).Value = "Poetry"
$ws.Cells.Item(1253, 5).Value = "Some artists write poetry"
$ws.Cells.Item(1253, 6).Value = "write poetry / sáng tác thơ"
$ws.Cells.Item(1253, 7).Value = "N"

$ws.Cells.Item(1254, 1).Value = "U42_23"
$ws.Cells.Item(1254, 2).Value = 42
$ws.Cells.Item(1254, 3).Value = "Phi thường"
$ws.Cells.Item(1254, 4).Value = "Extraordinary"
$ws.Cells.Item(1254, 5).Value = "My hero is an extraordinary person"
$ws.Cells.Item(1254, 6).Value = "An extraordinary person / một người phi thường"
$ws.Cells.Item(1254, 7).Value = "Adj"

$ws.Cells.Item(1255, 1).Value = "U42_24"
$ws.Cells.Item(1255, 2).Value = 42
$ws.Cells.Item(1255, 3).Value = "Thế kỷ"
$ws.Cells.Item(1255, 4).Value = "Century"
$ws.Cells.Item(1255, 5).Value = "This device was invented a century ago"
$ws.Cells.Item(1255, 6).Value = "a century ago / một thế kỷ trước"
$ws.Cells.Item(1255, 7).Value = "N"

$ws.Cells.Item(1256, 1).Value = "U42_25"
$ws.Cells.Item(1256, 2).Value = 42
$ws.Cells.Item(1256, 3).Value = "Chủ đề"
$ws.Cells.Item(1256, 4).Value = "Theme"
$ws.Cells.Item(1256, 5).Value = "The theme of the party is halloween"
$ws.Cells.Item(1256, 6).Value = "theme of something / chủ đề của một cái gì đó"
$ws.Cells.Item(1256, 7).Value = "N"

$ws.Cells.Item(1257, 1).Value = "U42_26"
$ws.Cells.Item(1257, 2).Value = 42
$ws.Cells.Item(1257, 3).Value = "Thiết kế"
$ws.Cells.Item(1257, 4).Value = "Design"
$ws.Cells.Item(1257, 5).Value = "It costs money to design a logo"
$ws.Cells.Item(1257, 6).Value = "design something"
$ws.Cells.Item(1257, 7).Value = "V"

$ws.Cells.Item(1258, 1).Value = "U42_27"
$ws.Cells.Item(1258, 2).Value = 42
$ws.Cells.Item(1258, 3).Value = "Biểu tượng"
$ws.Cells.Item(1258, 4).Value = "Symbol"
$ws.Cells.Item(1258, 5).Value = "Green is the symbol of nature"
$ws.Cells.Item(1258, 6).Value = "Symbol of something"
$ws.Cells.Item(1258, 7).Value = "N"

$ws.Cells.Item(1259, 1).Value = "U42_28"
$ws.Cells.Item(1259, 2).Value = 42
$ws.Cells.Item(1259, 3).Value = "Con rối"
$ws.Cells.Item(1259, 4).Value = "Puppet"
$ws.Cells.Item(1259, 5).Value = "The party will have a puppet show"
$ws.Cells.Item(1259, 6).Value = "a puppet show / một mà biểu diễn múa rối"
$ws.Cells.Item(1259, 7).Value = "N"

$ws.Cells.Item(1260, 1).Value = "U42_29"
$ws.Cells.Item(1260, 2).Value = 42
$ws.Cells.Item(1260, 3).Value = "Tài năng"
$ws.Cells.Item(1260, 4).Value = "Talented"
$ws.Cells.Item(1260, 5).Value = "In basketball, he is a talented player"
$ws.Cells.Item(1260, 6).Value = "a talented player"
$ws.Cells.Item(1260, 7).Value = "Adj"

$ws.Cells.Item(1261, 1).Value = "U42_30"
$ws.Cells.Item(1261, 2).Value = 42
$ws.Cells.Item(1261, 3).Value = "Tâm hồn"
$ws.Cells.Item(1261, 4).Value = "Soul"
$ws.Cells.Item(1261, 5).Value = "You have a heart of gold (trái tim nhân hậu) and a beautiful soul"
$ws.Cells.Item(1261, 6).Value = "a beautiful soul / tâm hồn đẹp"
$ws.Cells.Item(1261, 7).Value = "N"

$ws.Range("E1251").Select()
